$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '261.16'
$ws.Range('E2').NumberFormat = "@"
$ws.Range('E2').Value = '-0.21%'
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '27.06'
$ws.Range('E3').NumberFormat = "@"
$ws.Range('E3').Value = '-0.67%'
$ws.Range('D4').NumberFormat = "@"
$ws.Range('D4').Value = '4.696'
$ws.Range('E4').NumberFormat = "@"
$ws.Range('E4').Value = '-0.44%'
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '0.06224'
$ws.Range('E5').NumberFormat = "@"
$ws.Range('E5').Value = '2.60%'
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '6.746'
$ws.Range('E6').NumberFormat = "@"
$ws.Range('E6').Value = '1.57%'
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.8526'
$ws.Range('E7').NumberFormat = "@"
$ws.Range('E7').Value = '-1.16%'
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.9113'
$ws.Range('E8').NumberFormat = "@"
$ws.Range('E8').Value = '-1.26%'
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.1399'
$ws.Range('E9').NumberFormat = "@"
$ws.Range('E9').Value = '-0.64%'
$ws.Range('B10').NumberFormat = "@"
$ws.Range('B10').Value = 'MandalaExchangeToken'
$ws.Range('C10').NumberFormat = "@"
$ws.Range('C10').Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '0.07090'
$ws.Range('E10').NumberFormat = "@"
$ws.Range('E10').Value = '-0.09%'
$ws.Range('B11').NumberFormat = "@"
$ws.Range('B11').Value = 'BitrueCoin'
$ws.Range('C11').NumberFormat = "@"
$ws.Range('C11').Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.03130'
$ws.Range('E11').NumberFormat = "@"
$ws.Range('E11').Value = '2.67%'
$ws.Range('B12').NumberFormat = "@"
$ws.Range('B12').Value = 'BitMartToken'
$ws.Range('C12').NumberFormat = "@"
$ws.Range('C12').Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '0.09052'
$ws.Range('E12').NumberFormat = "@"
$ws.Range('E12').Value = '-0.53%'
$ws.Range('B13').NumberFormat = "@"
$ws.Range('B13').Value = 'BitForexToken'
$ws.Range('C13').NumberFormat = "@"
$ws.Range('C13').Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '0.001535'
$ws.Range('E13').NumberFormat = "@"
$ws.Range('E13').Value = '0.29%'
$ws.Range('B14').NumberFormat = "@"
$ws.Range('B14').Value = 'One'
$ws.Range('C14').NumberFormat = "@"
$ws.Range('C14').Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '0.0006170'
$ws.Range('E14').NumberFormat = "@"
$ws.Range('E14').Value = '1.56%'
$ws.Range('B15').NumberFormat = "@"
$ws.Range('B15').Value = 'TigerCash'
$ws.Range('C15').NumberFormat = "@"
$ws.Range('C15').Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '0.006131'
$ws.Range('E15').NumberFormat = "@"
$ws.Range('E15').Value = '1.37%'
$ws.Range('B16').NumberFormat = "@"
$ws.Range('B16').Value = 'LEO'
$ws.Range('C16').NumberFormat = "@"
$ws.Range('C16').Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '3.451'
$ws.Range('E16').NumberFormat = "@"
$ws.Range('E16').Value = '0.17%'
$ws.Range('B17').NumberFormat = "@"
$ws.Range('B17').Value = 'GateToken'
$ws.Range('C17').NumberFormat = "@"
$ws.Range('C17').Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '3.175'
$ws.Range('E17').NumberFormat = "@"
$ws.Range('E17').Value = '0.17%'
$ws.Range('B18').NumberFormat = "@"
$ws.Range('B18').Value = 'BTSEToken'
$ws.Range('C18').NumberFormat = "@"
$ws.Range('C18').Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '2.166'
$ws.Range('E18').NumberFormat = "@"
$ws.Range('E18').Value = '-0.49%'
$ws.Range('B19').NumberFormat = "@"
$ws.Range('B19').Value = 'BitpandaEcosystemToken'
$ws.Range('C19').NumberFormat = "@"
$ws.Range('C19').Value = 'https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best'
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '0.3106'
$ws.Range('E19').NumberFormat = "@"
$ws.Range('E19').Value = '-0.68%'
$ws.Range('B20').NumberFormat = "@"
$ws.Range('B20').Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range('C20').NumberFormat = "@"
$ws.Range('C20').Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '0.04718'
$ws.Range('E20').NumberFormat = "@"
$ws.Range('E20').Value = '-7.78%'
$ws.Range('E21').NumberFormat = "@"
$ws.Range('E21').Value = '0.94%'
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '4.100'
$ws.Range('E22').NumberFormat = "@"
$ws.Range('E22').Value = '0.29%'
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '0.04251'
$ws.Range('E23').NumberFormat = "@"
$ws.Range('E23').Value = '-0.37%'
$ws.Range('E24').NumberFormat = "@"
$ws.Range('E24').Value = '-0.38%'
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '0.004091'
$ws.Range('E25').NumberFormat = "@"
$ws.Range('E25').Value = '4.59%'
$ws.Range('E26').NumberFormat = "@"
$ws.Range('E26').Value = '-0.01%'
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '0.0001639'
$ws.Range('E27').NumberFormat = "@"
$ws.Range('E27').Value = '4.34%'
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '0.03903'
$ws.Range('E40').NumberFormat = "@"
$ws.Range('E40').Value = '0.57%'
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '0.004129'
$ws.Range('E42').NumberFormat = "@"
$ws.Range('E42').Value = '-0.09%'
$ws.Range('E43').NumberFormat = "@"
$ws.Range('E43').Value = '0.97%'
$ws.Range('E44').NumberFormat = "@"
$ws.Range('E44').Value = '-7.52%'
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '0.00005111'
$ws.Range('E45').NumberFormat = "@"
$ws.Range('E45').Value = '-3.82%'
$ws.Range('E46').NumberFormat = "@"
$ws.Range('E46').Value = '-0.02%'
$ws.Range('E47').NumberFormat = "@"
$ws.Range('E47').Value = '-37.68%'
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '0.06945'
$ws.Range('E48').NumberFormat = "@"
$ws.Range('E48').Value = '-48.67%'
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '0.00002100'
$ws.Range('E49').NumberFormat = "@"
$ws.Range('E49').Value = '-0.02%'
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '0.0002000'
$ws.Range('E50').NumberFormat = "@"
$ws.Range('E50').Value = '-0.02%'
